$wb = $excel.ActiveWorkbook

# Sheet 1 = 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(3, 6).Value = 5148
$ws1.Cells.Item(4, 6).Value = 10
$ws1.Cells.Item(5, 6).Value = 7447
$ws1.Cells.Item(11, 6).Value = 28
$ws1.Cells.Item(12, 6).Value = 4317
$ws1.Cells.Item(13, 6).Value = 1760
$ws1.Cells.Item(15, 6).Value = 106
$ws1.Cells.Item(16, 6).Value = 2919
$ws1.Cells.Item(17, 6).Value = 580
$ws1.Cells.Item(20, 6).Value = 501
$ws1.Cells.Item(21, 6).Value = 436
$ws1.Cells.Item(22, 6).Value = 457
$ws1.Cells.Item(23, 6).Value = 307
$ws1.Cells.Item(24, 6).Value = 100
$ws1.Cells.Item(25, 6).Value = 1694
$ws1.Cells.Item(26, 6).Value = 1185
$ws1.Cells.Item(28, 6).Value = 1378
$ws1.Cells.Item(30, 6).Value = 579
$ws1.Cells.Item(31, 6).Value = 27
$ws1.Cells.Item(36, 6).Value = 65
$ws1.Cells.Item(37, 6).Value = 2889
$ws1.Cells.Item(38, 6).Value = 705
$ws1.Cells.Item(39, 6).Value = 20
$ws1.Cells.Item(40, 6).Value = 72
$ws1.Cells.Item(42, 6).Value = 28

# Sheet 2 = 演出 (Performances)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(2, 6).Value = 11

# Sheet 4 = 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(3, 6).Value = 5148
$ws4.Cells.Item(4, 6).Value = 10
$ws4.Cells.Item(5, 6).Value = 7447
$ws4.Cells.Item(11, 6).Value = 28
$ws4.Cells.Item(12, 6).Value = 4317
$ws4.Cells.Item(13, 6).Value = 1760
$ws4.Cells.Item(15, 6).Value = 106
$ws4.Cells.Item(16, 6).Value = 2919
$ws4.Cells.Item(17, 6).Value = 580
$ws4.Cells.Item(20, 6).Value = 501
$ws4.Cells.Item(21, 6).Value = 436
$ws4.Cells.Item(22, 6).Value = 457
$ws4.Cells.Item(23, 6).Value = 11
$ws4.Cells.Item(24, 6).Value = 307
$ws4.Cells.Item(25, 6).Value = 100
$ws4.Cells.Item(26, 6).Value = 1694
$ws4.Cells.Item(27, 6).Value = 1185
$ws4.Cells.Item(29, 6).Value = 1378
$ws4.Cells.Item(31, 6).Value = 579
$ws4.Cells.Item(32, 6).Value = 27
$ws4.Cells.Item(37, 6).Value = 65
$ws4.Cells.Item(38, 6).Value = 2889
$ws4.Cells.Item(40, 6).Value = 705
$ws4.Cells.Item(41, 6).Value = 20
$ws4.Cells.Item(42, 6).Value = 72
$ws4.Cells.Item(44, 6).Value = 28
